$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before B; this shifts old column B (Valor data) to column C
$ws.Columns("B").Insert()

# Header row
$ws.Range("B1").Value = "Variável"
$ws.Range("D1").Value = "Colocação"

# Column B: variable label for each data row (2-10)
$ws.Range("B2:B10").Value = "Variação 2021/2020"

# Column D: ranking labels for rows 2-8 only
$ws.Range("D2").Value = "1º"
$ws.Range("D3").Value = "2º"
$ws.Range("D4").Value = "3º"
$ws.Range("D5").Value = "4º"
$ws.Range("D6").Value = "5º"
$ws.Range("D7").Value = "6º"
$ws.Range("D8").Value = "13º"
